$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update input values (column B) ---
$ws.Range("B2").Value = 0
$ws.Range("B5").Value = 4.4000000000000004
$ws.Range("B6").Value = 4.4000000000000004
$ws.Range("B8").Value = 207.7
$ws.Range("B9").Value = 5010
$ws.Range("B10").Value = 152
$ws.Range("B15").Value = 4.2
$ws.Range("B16").Value = 4.87
$ws.Range("B18").Value = 3600

# --- Update formulas ---
# G5 used to hold the static note "See D5"; now it mirrors E5 via formula
$ws.Range("G5").Formula = "=E5"

# H5 used to hold the same static note "See D5"; now it holds the label "E5"
$ws.Range("H5").Value = "E5"

# G7 used to hold the static note "10% of G6"; now it calculates 10% of G6
$ws.Range("G7").Formula = "=IMPRODUCT(G6,0.1)"

# C8 / C9 formulas had their multipliers swapped
$ws.Range("C8").Formula = "=IMPRODUCT(B8, 3600)"
$ws.Range("C9").Formula = "=IMPRODUCT(B9,60)"

# New E8 / E9 formulas
$ws.Range("E8").Formula = "=IMDIV(C8, B10)"
$ws.Range("E9").Formula = "=IMDIV(C9, B10)"

# --- Update selection / active cell ---
[void]$ws.Range("B11").Select()
